$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Update the benchmark value for I4 (formulas in J4, J12, J13 recalc automatically)
$ws.Range("I4").Value = 157.14500000000001

# Move the active selection to I4 (matches saved view state in the diff)
$ws.Activate()
$ws.Range("I4").Select()
